$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ST_STAT_01 - "Loc thong ke theo ngay" ---
# Cac Buoc (C2): append "(Dung JS)" to the date line
$ws.Range("C2").Value = "1. Menu -> Báo cáo chi tiết`n2. Nhập ngày 2025-10-01 - 2025-12-31 (Dùng JS)`n3. Bấm Xem kết quả"

# Ket Qua Mong Doi (E2): clarify expected result
$ws.Range("E2").Value = "Hiển thị doanh thu và danh sách sản phẩm bán chạy"

# Ket Qua Thuc Te (F2): updated observed numbers
$ws.Range("F2").Value = "Doanh thu hiển thị: 2,500,000 đ | Số dòng SP: 1"

# --- Row 3: ST_STAT_02 - "Xuat bao cao Excel" ---
# Cac Buoc (C3): rewritten steps
$ws.Range("C3").Value = "1. Điều hướng lại`n2. Set lại ngày (2025-10-01 - 2025-12-31)`n3. Bấm nút 'Xuất file Excel'"

# Ket Qua Mong Doi (E3): rewritten expected result
$ws.Range("E3").Value = "Server xử lý request và trình duyệt tải file (.xlsx), không báo lỗi Server (500)"

# --- Column widths re-fitted to the new content ---
$ws.Columns.Item(3).ColumnWidth = 43.333333333333336
$ws.Columns.Item(5).ColumnWidth = 67.333333333333336
